# Apply the edits described by the commit:
#  - rename sheet "Sheet2" -> "Sheet1"
#  - add new "periode" column at the very front (new column A),
#    shifting every existing column one to the right
#  - add a new "aset_likuid_tidak_menghasilkan" column right after
#    "aset tidak menghasilkan" (E7 in PEARLS terms)
#  - drop the old trailing "periode" / "tgl buat" columns and replace them
#    with a single new "tanggal buat" column (L2 in PEARLS terms) at the end
#  - move the active selection to C11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Sheet1"

# Remove the old trailing "periode" and "tgl buat" header columns (AK:AL)
$ws.Columns("AK:AL").Delete()

# Insert a brand-new first column and give it the "periode" header
$ws.Columns("A:A").Insert()
$ws.Range("A1").Value = "periode"

# Insert the new "aset_likuid_tidak_menghasilkan" column right after
# "aset tidak menghasilkan" (now column J) -> new column K
$ws.Columns("K:K").Insert()
$ws.Range("K1").Value = "aset_likuid_tidak_menghasilkan"

# Append the new "tanggal buat" header as the new last column
$lastCol = $ws.UsedRange.Columns.Count
$ws.Cells(1, $lastCol + 1).Value = "tanggal buat"

# Restore the selection shown in the saved workbook
$ws.Range("C11").Select()
